# "petit bug de réglé" -- a bug in the trial-generation script was fixed and
# rows 2 and 3 (the first two generated trials) were re-rolled with corrected
# stimulus/distractor pairings and new reaction-time (TR) values.
#
# Columns: A=Stimulus, B=Déterminant_Mot, C=Nom_Mot, D=Déterminant_image,
#          E=Nom_Image, F=TR, G=Lettre, H=Congruence, I=Erreur

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new Déterminant/Nom pairing + new measured reaction time -----
$ws.Range("B2").Value = "Les"
$ws.Range("C2").Value = "oignons"
$ws.Range("D2").Value = "Un"
$ws.Range("E2").Value = "oignon"
$ws.Range("F2").Value = 0.56456759991124272

# --- Row 3: new Déterminant/Nom pairing + new measured reaction time -----
$ws.Range("B3").Value = "Les"
$ws.Range("C3").Value = "citrouilles"
$ws.Range("D3").Value = "Les"
$ws.Range("E3").Value = "citrouilles"
$ws.Range("F3").Value = 0.58663499995600432

# The regeneration run also re-wrote the formatting of the header row and of
# the two re-rolled rows (fresh style entries in the workbook's style
# table). Re-apply the (unchanged) text format to those same cells so they
# pick up fresh style entries too.
$ws.Range("A1:I1").NumberFormat = "@"
$ws.Range("A2:E2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("A3:E3").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
